$wb = $excel.ActiveWorkbook

# Sheets: test_suite, addCustomer, OpenAccount
$wsSuite = $wb.Worksheets.Item("test_suite")
$wsAddCustomer = $wb.Worksheets.Item("addCustomer")
$wsOpenAccount = $wb.Worksheets.Item("OpenAccount")

# --- test_suite: fix AddCustomer's runmode from N to Y ---
$wsSuite.Range("B3").Value = "Y"
$wsSuite.Activate()
$wsSuite.Range("B3").Select() | Out-Null

# --- addCustomer: add "runmode" column (D) set to Y for all rows ---
$wsAddCustomer.Range("D1").Value = "runmode"
$wsAddCustomer.Range("D2").Value = "Y"
$wsAddCustomer.Range("D3").Value = "Y"
$wsAddCustomer.Range("D4").Value = "Y"
$wsAddCustomer.Range("D5").Value = "Y"
$wsAddCustomer.Activate()
$wsAddCustomer.Range("D5").Select() | Out-Null

# --- OpenAccount: add "runmode" column (C) set to Y ---
$wsOpenAccount.Range("C1").Value = "runmode"
$wsOpenAccount.Range("C2").Value = "Y"
$wsOpenAccount.Activate()
$wsOpenAccount.Range("C2").Select() | Out-Null

# Re-activate addCustomer sheet as the active tab (tabSelected)
$wsAddCustomer.Activate()
